# Fix login-module bug config: unify the AnimaState ("NormalSkillX") values
# for the per-hero "Normal" skill rows (12-19) so they all reference the
# same shared string as row 11, and rename the "SpecialAttackN" skill
# identifiers to "SpecialSkillN" everywhere they are used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Make C12:C19 match C11's value ("NormalSkill1,NormalSkill2,NormalSkill3")
$normalSkillValue = $ws.Range("C11").Value2
for ($row = 12; $row -le 19; $row++) {
    $ws.Range("C" + $row).Value2 = $normalSkillValue
}

# 2. Rename "SpecialAttack1" -> "SpecialSkill1" and "SpecialAttack2" -> "SpecialSkill2"
#    across the whole used range (column C on the "Attack"/"THUMP" rows).
$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $val = $cell.Value2
    if ($val -eq "SpecialAttack1") {
        $cell.Value2 = "SpecialSkill1"
    } elseif ($val -eq "SpecialAttack2") {
        $cell.Value2 = "SpecialSkill2"
    }
}

# 3. Update the view state: scroll so row 50 is the frozen pane's top-left
#    cell, and move/replace the active selection to C62.
$ws.Activate()
$ws.Range("A50").Select()
$ws.Application.ActiveWindow.ScrollRow = 50
$ws.Range("C62").Select()
